$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix markdown rendering: replace literal "<br>" line-break markers in the
# --- GCS coverage-tier descriptions (column E) with real line breaks so the
# --- text renders correctly instead of showing literal "<br>" tags.

$ws.Range("E26:E37").Value = "Supports the GCS if coverage is **Low** `nOther members: Global South + EU `n(25-33% of world emissions)"
$ws.Range("E38:E49").Value = "Supports the GCS if coverage is **Mid** `nGlobal South + China `n(56% of world emissions)"
$ws.Range("E50:E61").Value = "Supports the GCS if coverage is **High** `nGlobal South + China + EU + various HICs `n(UK, Japan, Korea, Canada...; 64-72% of emissions)"
$ws.Range("E62:E73").Value = "Supports the GCS if coverage is **High**, **color** variant `nGlobal South + China + EU + various HICs `n+ Distributive effects shown using colors on world map"

# --- Updated model-output statistics (mean / CI_low / CI_high) for the
# --- "<b>All</b>" and "Russia" rows of each 12-row variable block.

$ws.Range("B2").Value = 67.7971034621066
$ws.Range("C2").Value = 66.4060461449477
$ws.Range("D2").Value = 69.1881607792655

$ws.Range("B12").Value = 73.8102295238027
$ws.Range("C12").Value = 69.9173212117554
$ws.Range("D12").Value = 77.70313783585

$ws.Range("B14").Value = 55.3794402736958
$ws.Range("C14").Value = 54.4904084313456
$ws.Range("D14").Value = 56.268472116046

$ws.Range("B24").Value = 49.0525173251567
$ws.Range("C24").Value = 46.6565838643414
$ws.Range("D24").Value = 51.448450785972

$ws.Range("B38").Value = 67.1068857857389
$ws.Range("C38").Value = 65.398148692157
$ws.Range("D38").Value = 68.8156228793208

$ws.Range("B48").Value = 63.180865855067
$ws.Range("C48").Value = 57.4064206273406
$ws.Range("D48").Value = 68.9553110827935

$ws.Range("B50").Value = 68.4640848889378
$ws.Range("C50").Value = 66.8324549522336
$ws.Range("D50").Value = 70.095714825642

$ws.Range("B60").Value = 59.8978761247204
$ws.Range("C60").Value = 53.8467358677426
$ws.Range("D60").Value = 65.9490163816982

$ws.Range("B62").Value = 61.8966560897
$ws.Range("C62").Value = 60.1646532127507
$ws.Range("D62").Value = 63.6286589666493

$ws.Range("B72").Value = 54.1378328703989
$ws.Range("C72").Value = 47.9978727513762
$ws.Range("D72").Value = 60.2777929894216
